$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Workbook window view (bookViews/workbookView) - best effort; the
# geometry is purely cosmetic (maximized-window rectangle) and has no
# effect on workbook data.
# ---------------------------------------------------------------------
$win = $wb.Windows.Item(1)
$win.Left = -38400
$win.Top = 460
$win.Width = 38400
$win.Height = 21060

# ---------------------------------------------------------------------
# Sheet references
# ---------------------------------------------------------------------
$wsFull = $wb.Worksheets.Item("Full Instrument")
$wsShield = $wb.Worksheets.Item("Shield Parts")
$wsToAdd = $wb.Worksheets.Item("To Add")

# ---------------------------------------------------------------------
# Shield Parts: D1 note
# ---------------------------------------------------------------------
$wsShield.Range("D1").Value = "Shipping Costs Not Included!"

# ---------------------------------------------------------------------
# Shield Parts: row 4 (new SHT21 / DEV-11114 part)
# ---------------------------------------------------------------------
$wsShield.Range("C4").Value = "Arduino Pro Mini"
$wsShield.Range("D4").Value = 1
$wsShield.Range("E4").Value = 9.9499999999999993
$wsShield.Range("G4").Value = 9.4499999999999993
$wsShield.Range("I4").Value = 8.9600000000000009
$wsShield.Range("K4").Value = "Sparkfun"
$wsShield.Range("L4").Value = "DEV-11114"
$wsShield.Range("M4").Value = "N/A"
$wsShield.Range("N4").Value = "https://www.sparkfun.com/products/11114"
$wsShield.Range("O4").Value = "3.3V/8MHz Version"

# ---------------------------------------------------------------------
# Full Instrument: row 5 (SHT21 Ebay source)
# ---------------------------------------------------------------------
$wsFull.Range("E5").Value = 11.2
$wsFull.Range("G5").Value = 11.2
$wsFull.Range("I5").Value = 11.2
$wsFull.Range("K5").Value = "Ebay"
$wsFull.Range("L5").Value = "N/A"
$wsFull.Range("N5").Value = "http://www.ebay.com/itm/High-Quality-Humidity-Sensor-SHT21-Breakout-Board-GY-21-/281654379137"
$wsFull.Range("O5").Value = "May need alternates. Not reliable source"

# ---------------------------------------------------------------------
# Shield Parts: row 9 (Break away headers)
# ---------------------------------------------------------------------
$wsShield.Range("C9").Value = "Break Away Headers - Straight"
$wsShield.Range("D9").Value = 1
$wsShield.Range("E9").Value = 1.5
$wsShield.Range("G9").Value = 1.5
$wsShield.Range("I9").Value = 1.43
$wsShield.Range("K9").Value = "Sparkfun"
$wsShield.Range("L9").Value = "PRT-00116"
$wsShield.Range("M9").Value = "40-Pin 0.1"" Through Hole Header"
$wsShield.Range("N9").Value = "https://www.sparkfun.com/products/116"
$wsShield.Range("O9").Value = "Comes in strips of 40, currently only need 34. Prices will be updated once precise count is established"

# ---------------------------------------------------------------------
# To Add: remove the "Arduino Pro Mini" line (now tracked on Shield Parts),
# update the "Headers" line, and append 3 more follow-up items.
# ---------------------------------------------------------------------
$wsToAdd.Range("A2").ClearContents()
$wsToAdd.Range("A3").Value = "Sensor Headers - Add to full instrument"
$wsToAdd.Range("A4").Value = "Jumper Wires - Add to full instrument"
$wsToAdd.Range("A5").Value = "Shield Resistors"
$wsToAdd.Range("A6").Value = "Shield LEDs"

# ---------------------------------------------------------------------
# Column widths - the newly-added long URLs/notes push these columns
# wider (best-fit recalculation). Set explicit character widths to
# match the post-edit best-fit result.
# ---------------------------------------------------------------------
$wsFull.Columns("N:N").ColumnWidth = 86.1640625
$wsFull.Columns("O:O").ColumnWidth = 34.33203125
$wsShield.Columns("M:M").ColumnWidth = 28
$wsShield.Columns("O:O").ColumnWidth = 83.33203125

# ---------------------------------------------------------------------
# Selections (cursor position when the file was last saved)
# ---------------------------------------------------------------------
[void]$wsFull.Range("C39").Select()
[void]$wsShield.Range("C28").Select()
[void]$wsToAdd.Range("A7").Select()
[void]$wsToAdd.Activate()

[void]$wb.Application.Calculate()
